# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right after "总计" (i.e. before "2022-Q2"),
#    built as a duplicate of "2022-Q2" (same headers/layout/formatting) with
#    its fund figures updated to the 2022-Q3 numbers.
# 2) Insert a new summary row into "总计" for "2022-Q3", pushing the existing
#    quarters down by one row and keeping the running index in column A
#    consistent.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate "2022-Q2" -> "2022-Q3", placed right after "总计"
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $summary)

$q3 = $wb.Worksheets.Item("总计").Next
$q3.Name = "2022-Q3"

# these columns hold text-formatted figures (e.g. "0.46", not the number
# 0.46) in the source data, so force text with a leading apostrophe rather
# than letting Excel auto-detect them as numbers
$q3.Range("D2").Value = "'0.46"
$q3.Range("E2").Value = "'76.66"
$q3.Range("F2").Value = "'3.96"
$q3.Range("G2").Value = "'0.0182"
$q3.Range("H2").Value = 10

$q3.Range("D3").Value = "'0.08"
$q3.Range("E3").Value = "'76.66"
$q3.Range("F3").Value = "'3.96"
$q3.Range("G3").Value = "'0.0032"
$q3.Range("H3").Value = 10

# ---------------------------------------------------------------------
# Step 2: insert the "2022-Q3" row into "总计", shifting older quarters down
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")

$lastRow = 7
for ($r = $lastRow; $r -ge 2; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dst, 4).Value = $ws.Cells.Item($r, 4).Value2
}

# the shift above leaves the newly extended last row without the column-A
# formatting the rest of the table has - copy it across before filling values
$ws.Range("A2").Copy()
$ws.Range("A" + ($lastRow + 1)).PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "2022-Q3"
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = 0.02

# keep the running index in column A (0,1,2,...) consistent after the insert
for ($r = 3; $r -le ($lastRow + 1); $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# restore the originally active tab (our sheet copy/rename above moved it)
$wb.Worksheets.Item("2020-Q4").Activate()
